$d = $word.ActiveDocument

$d.Content.Find.Execute("16:30h - 18:00h - Introdução à espectrometria de massas", $true, $false, $false, $false, $false,
                         $true, 1, $false, "16:00h - 18:00h - Introdução à espectrometria de massas", 2)

$d.Content.Find.Execute("16:30h - 18:00h - Proteoma Quantitativo: aplicações e estratégias", $true, $false, $false, $false, $false,
                         $true, 1, $false, "16:00h - 18:00h - Proteoma Quantitativo: aplicações e estratégias", 2)
